$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.609.21"
$ws.Range("E2").Value = "  -1.54%  "

$ws.Range("D3").Value = "3.523.94"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "609.28"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").Value = "183.56"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -1.76%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "0.211"
$ws.Range("E9").Value = "  +4.07%  "

$ws.Range("D10").Value = "0.638"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("D11").Value = "53.27"
$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "9.39"
$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("D14").Value = "4.088.47"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").Value = "587.26"
$ws.Range("E15").Value = "  +4.43%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.690.65"
$ws.Range("E16").Value = "  -1.39%  "

$ws.Range("D17").Value = "3.527.37"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "18.80"
$ws.Range("E19").Value = "  -4.48%  "

$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").Value = "0.985"
$ws.Range("E21").Value = "  -3.14%  "

$ws.Range("D22").Value = "17.41"
$ws.Range("E22").Value = "  -3.08%  "

$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "96.62"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "4.82"
$ws.Range("E25").Value = "  -1.95%  "

$ws.Range("D26").Value = "2.95"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -6.24%  "

$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  +3.82%  "

$ws.Range("D29").Value = "31.91"
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -4.53%  "

$ws.Range("D31").Value = "12.08"
$ws.Range("E31").Value = "  -3.24%  "

$ws.Range("D33").Value = "63.20"
$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("D35").Value = "3.58"
$ws.Range("E35").Value = "  +16.34%  "

$ws.Range("D36").Value = "531.34"
$ws.Range("E36").Value = "  -5.86%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.399"
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("D39").Value = "36.88"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("D40").Value = "3.532.75"
$ws.Range("E40").Value = "  +5.33%  "

$ws.Range("D41").Value = "0.0₃0774"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("E42").Value = "  +3.82%  "

$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("D44").Value = "0.0453"
$ws.Range("E44").Value = "  +1.36%  "

$ws.Range("D45").Value = "2.91"
$ws.Range("E45").Value = "  -2.33%  "

$ws.Range("D46").Value = "3.38"
$ws.Range("E46").Value = "  -4.61%  "

$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("D48").Value = "9.08"
$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("E50").Value = "  -6.14%  "

$ws.Range("D51").Value = "134.38"
$ws.Range("E51").Value = "  -2.42%  "
